$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.127.56"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.667.99"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'210.78"
$ws.Range("E5").Value = "  -2.42%  "
$ws.Range("D6").Value = "'0.5221"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "'0.2622"
$ws.Range("E8").Value = "  -2.46%  "
$ws.Range("D9").Value = "'0.06324"
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("D10").Value = "'21.17"
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("D11").Value = "'0.07552"
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").Value = "1.667.03"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").Value = "'4.426"
$ws.Range("E13").Value = "  -2.16%  "
$ws.Range("D14").Value = "'0.5480"
$ws.Range("E14").Value = "  -4.93%  "
$ws.Range("D15").Value = "'0.000008017"
$ws.Range("E15").Value = "  -2.80%  "
$ws.Range("D16").Value = "'66.38"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "26.162.20"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "'4.744"
$ws.Range("E19").Value = "  -2.59%  "
$ws.Range("D20").Value = "'187.34"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").Value = "'10.29"
$ws.Range("E21").Value = "  -4.37%  "
$ws.Range("D22").Value = "'6.232"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "'149.51"
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("E25").Value = "  -1.81%  "
$ws.Range("D26").Value = "'7.467"
$ws.Range("E26").Value = "  -3.31%  "
$ws.Range("D27").Value = "'15.80"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("D29").Value = "'1.353"
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("D31").Value = "'3.523"
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("D32").Value = "'3.414"
$ws.Range("E32").Value = "  -4.37%  "
$ws.Range("D33").Value = "'1.646"
$ws.Range("E33").Value = "  -2.17%  "
$ws.Range("D34").Value = "'1.004"
$ws.Range("E34").Value = "  -1.63%  "
$ws.Range("D35").Value = "'0.6028"
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("D36").Value = "'2.397"
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("D38").Value = "1.116.38"
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("D39").Value = "'6.077"
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("D40").Value = "'0.01611"
$ws.Range("E40").Value = "  -1.69%  "
$ws.Range("D41").Value = "'0.8639"
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").Value = "'100.48"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "1.822.06"
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("D45").Value = "'0.00000000111"
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("D46").Value = "'55.41"
$ws.Range("E46").Value = "  -3.71%  "
$ws.Range("D47").Value = "'0.9996"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "'8.065"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").Value = "'0.05234"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").Value = "'5.916"
$ws.Range("E51").Value = "  -1.70%  "
